$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4968.864
$ws.Range("I33").Value = 6414.625
$ws.Range("J33").Value = 1113.5
$ws.Range("K33").Value = 6414.625
$ws.Range("L33").Value = 1113.5
$ws.Range("M33").Value = -6185.625
$ws.Range("N33").Value = -1571.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 37037760
$ws.Range("I101").Value = 41666980
$ws.Range("K101").Value = 125000940
$ws.Range("M101").Value = -124999318

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 35150
$ws.Range("J117").Value = 35150
$ws.Range("L117").Value = 35150
$ws.Range("N117").Value = -44328

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3942.3215
$ws.Range("I138").Value = 1410.8667
$ws.Range("J138").Value = 6863.231
$ws.Range("K138").Value = 4232.6001
$ws.Range("L138").Value = 20589.693
$ws.Range("M138").Value = 907.3999000000003
$ws.Range("N138").Value = -30869.693

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1469.85
$ws.Range("I2").Value = 1239.4
$ws.Range("J2").Value = 2161.2
$ws.Range("K2").Value = 1239.4
$ws.Range("L2").Value = 2161.2
$ws.Range("M2").Value = -1126.4
$ws.Range("N2").Value = -2387.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1469.85
$ws.Range("I116").Value = 1239.4
$ws.Range("J116").Value = 2161.2
$ws.Range("K116").Value = 1239.4
$ws.Range("L116").Value = 2161.2
$ws.Range("M116").Value = 1054.6
$ws.Range("N116").Value = -6749.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1469.85
$ws.Range("I3").Value = 1239.4
$ws.Range("J3").Value = 2161.2
$ws.Range("K3").Value = 1239.4
$ws.Range("L3").Value = 2161.2
$ws.Range("M3").Value = -1125.4
$ws.Range("N3").Value = -2389.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 2071.3333
$ws.Range("I25").Value = 1707
$ws.Range("J25").Value = 2800
$ws.Range("K25").Value = 1707
$ws.Range("L25").Value = 2800
$ws.Range("M25").Value = -1472
$ws.Range("N25").Value = -3270

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2305.65
$ws.Range("I86").Value = 2460.25
$ws.Range("J86").Value = 2073.75
$ws.Range("K86").Value = 2460.25
$ws.Range("L86").Value = 2073.75
$ws.Range("M86").Value = -1337.25
$ws.Range("N86").Value = -4319.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2305.65
$ws.Range("I89").Value = 2460.25
$ws.Range("J89").Value = 2073.75
$ws.Range("K89").Value = 12301.25
$ws.Range("L89").Value = 10368.75
$ws.Range("M89").Value = -6685.25
$ws.Range("N89").Value = -21600.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 12401
$ws.Range("J92").Value = 12401
$ws.Range("L92").Value = 12401
$ws.Range("N92").Value = -17393

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1511.9354
$ws.Range("I99").Value = 874.61536
$ws.Range("K99").Value = 874.61536
$ws.Range("M99").Value = 623.38464

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3387.5
$ws.Range("I105").Value = 3157.1428
$ws.Range("K105").Value = 3157.1428
$ws.Range("M105").Value = -1410.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 933.25
$ws.Range("I22").Value = 1100.9231
$ws.Range("J22").Value = 206.66667
$ws.Range("K22").Value = 1100.9231
$ws.Range("L22").Value = 206.66667
$ws.Range("M22").Value = -750.9231
$ws.Range("N22").Value = -906.6666700000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 56426.223
$ws.Range("I122").Value = 77815.62
$ws.Range("J122").Value = 813.8
$ws.Range("K122").Value = 233446.86
$ws.Range("L122").Value = 2441.4
$ws.Range("M122").Value = -230996.86
$ws.Range("N122").Value = -7341.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 10546
$ws.Range("I6").Value = 15504
$ws.Range("J6").Value = 630
$ws.Range("K6").Value = 15504
$ws.Range("L6").Value = 630
$ws.Range("M6").Value = -15391
$ws.Range("N6").Value = -856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H16").Value = 10546
$ws.Range("I16").Value = 15504
$ws.Range("J16").Value = 630
$ws.Range("K16").Value = 15504
$ws.Range("L16").Value = 630
$ws.Range("M16").Value = -15254
$ws.Range("N16").Value = -1130

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1821.2858
$ws.Range("I126").Value = 1456.1428
$ws.Range("J126").Value = 2369
$ws.Range("K126").Value = 4368.428400000001
$ws.Range("L126").Value = 7107
$ws.Range("M126").Value = -1898.428400000001
$ws.Range("N126").Value = -12047

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3435.5908
$ws.Range("I132").Value = 2262.1667
$ws.Range("J132").Value = 4843.7
$ws.Range("K132").Value = 6786.500100000001
$ws.Range("L132").Value = 14531.1
$ws.Range("M132").Value = -4256.500100000001
$ws.Range("N132").Value = -19591.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1276.381
$ws.Range("I22").Value = 943.875
$ws.Range("J22").Value = 2340.4
$ws.Range("K22").Value = 943.875
$ws.Range("L22").Value = 2340.4
$ws.Range("M22").Value = -648.875
$ws.Range("N22").Value = -2930.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1276.381
$ws.Range("I27").Value = 943.875
$ws.Range("J27").Value = 2340.4
$ws.Range("K27").Value = 943.875
$ws.Range("L27").Value = 2340.4
$ws.Range("M27").Value = -836.875
$ws.Range("N27").Value = -2554.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 13917650
$ws.Range("I46").Value = 33400570
$ws.Range("J46").Value = 1279
$ws.Range("K46").Value = 33400570
$ws.Range("L46").Value = 1279
$ws.Range("M46").Value = -33400382
$ws.Range("N46").Value = -1655

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 209.23529
$ws.Range("I55").Value = 250.58333
$ws.Range("J55").Value = 110
$ws.Range("K55").Value = 250.58333
$ws.Range("L55").Value = 110
$ws.Range("M55").Value = -77.58332999999999
$ws.Range("N55").Value = -456

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2008
$ws.Range("I61").Value = 1415
$ws.Range("J61").Value = 2601
$ws.Range("K61").Value = 1415
$ws.Range("L61").Value = 2601
$ws.Range("M61").Value = -1213
$ws.Range("N61").Value = -3005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 12011.889
$ws.Range("I93").Value = 15129
$ws.Range("K93").Value = 15129
$ws.Range("M93").Value = -13881

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2008
$ws.Range("I113").Value = 1415
$ws.Range("J113").Value = 2601
$ws.Range("K113").Value = 1415
$ws.Range("L113").Value = 2601
$ws.Range("M113").Value = 755
$ws.Range("N113").Value = -6941

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3082.682
$ws.Range("I122").Value = 2187.75
$ws.Range("J122").Value = 4156.6
$ws.Range("K122").Value = 6563.25
$ws.Range("L122").Value = 12469.8
$ws.Range("M122").Value = -4113.25
$ws.Range("N122").Value = -17369.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 25433.8
$ws.Range("J133").Value = 25433.8
$ws.Range("L133").Value = 25433.8
$ws.Range("N133").Value = -30493.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 70011.336
$ws.Range("I21").Value = 70000
$ws.Range("K21").Value = 70000
$ws.Range("M21").Value = -69765

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 70011.336
$ws.Range("I35").Value = 70000
$ws.Range("K35").Value = 70000
$ws.Range("M35").Value = -69710

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5298

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 940.3
$ws.Range("I81").Value = 900.1667
$ws.Range("J81").Value = 1000.5
$ws.Range("K81").Value = 1800.3334
$ws.Range("L81").Value = 2001
$ws.Range("M81").Value = -739.3334
$ws.Range("N81").Value = -4123

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 940.3
$ws.Range("I84").Value = 900.1667
$ws.Range("J84").Value = 1000.5
$ws.Range("K84").Value = 9001.666999999999
$ws.Range("L84").Value = 10005
$ws.Range("M84").Value = -3697.666999999999
$ws.Range("N84").Value = -20613
